# Add a new "Title and Content" slide after the existing title slide,
# containing the Blackboard recording link (commit: "inclusion of Blackboard Link").

$p = $ppt.ActivePresentation

# "Title and Content" is the 2nd layout on the slide master.
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$s = $p.Slides.AddSlide(2, $layout)

# Title placeholder
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Lesson Video Link"
$title.TextFrame.TextRange.LanguageID = "en-GB"

# Content placeholder with the Blackboard Collaborate recording URL
$content = $s.Shapes.Item(2)
$content.TextFrame.TextRange.Text = "https://eu-lti.bbcollab.com/recording/3eba6015ec1b4492b52a514fc9174600"
$content.TextFrame.TextRange.LanguageID = "en-GB"

# Resize the content placeholder to hug the single line of link text.
$content.Left = 66
$content.Top = 143.75
$content.Width = 828
$content.Height = 92.06654
